# Update cryptocurrency price/volume data cells to reflect latest scrape
# (GitHub Actions symbol-list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.01%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.19%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.121"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.13%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07624"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.37%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-0.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.465"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.88%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9042"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.30%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1119"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'12.67%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1781"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.15%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-3.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.35%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001259"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.11%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005716"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.16%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'4.252"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.72%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.671"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.13%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2792"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.98%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04072"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001248"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'4.28%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004116"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.58%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.08%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003749"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02380"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.18%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05180"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.57%"
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'-1.90%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.81%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007055"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'10.02%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.001952"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.06%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007935"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-9.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-7.46%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'6.43%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.03175"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'814.65%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.13%"
$ws.Range("E51").Style = "Normal"
